# Carta_Incumplimiento_Informe.docx edit script
# Applies the changes described by the commit "Actualizar etiqueta [NÚMERO] en templates"

$d = $word.ActiveDocument

# --- helpers -----------------------------------------------------------
# Find `searchText` starting at/after `startPos` and before `endPos`.
# Returns the matched Range (its Start/End reflect the match), or $null.
function Find-From($startPos, $endPos, $searchText) {
    $scope = $d.Range($startPos, $endPos)
    $f = $scope.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
    if ($f) {
        return $scope
    }
    return $null
}

# Force Word to split the run(s) under `rng` into an isolated run with
# identical formatting (no visible side effect) by toggling a Font
# property on and back off. Word (and this host) re-materialises the
# run boundaries to exactly match `rng` when it writes the property.
function Isolate-Run($rng) {
    $rng.Font.Bold = 1
    $rng.Font.Bold = 0
}

# Walk a list of contiguous substrings inside [paraStart, paraEnd),
# isolating each one into its own run (so later formatting/content
# edits land on a single clean run instead of spilling onto siblings).
function Split-Sequence($paraStart, $paraEnd, $targets) {
    $pos = $paraStart
    foreach ($t in $targets) {
        $scope = Find-From $pos $paraEnd $t
        if ($scope -ne $null) {
            Isolate-Run $scope
            $pos = $scope.End
        }
    }
}

function Get-ParagraphRange($likeText) {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -like $likeText) {
            return $p.Range
        }
    }
    return $null
}

# --- Edit 1: "], entre InnovaChile y " -> split off "InnovaChile" -----
$para1 = Get-ParagraphRange "*entre InnovaChile y*"
Split-Sequence $para1.Start $para1.End @("], entre ", "InnovaChile", " y ")

# --- Edit 2: "N°[número]" -> "N° [NÚMERO]", then split into runs -----
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute("N°[número]", $true, $false, $false, $false, $false, $true, 1, $false, "N° [NÚMERO]", 2) | Out-Null

$para2 = Get-ParagraphRange "*aprobado por Resoluci*"
Split-Sequence $para2.Start $para2.End @("(E) ", "N°", " ", "[NÚMERO]", ", ", "de ")

# --- Edit 3: highlight "el día [día] de [mes] de [año]." ---------------
$find3 = $d.Content.Find
$find3.ClearFormatting()
$find3.Replacement.ClearFormatting()
$find3.Replacement.Highlight = $true
$find3.Execute("el día [día] de [mes] de [año].", $true, $false, $false, $false, $false, $true, 1, $false, `
    "el día [día] de [mes] de [año].", 2) | Out-Null

# --- Edit 4: "la garantía de fiel cumplimiento" -> "las garantías asociadas al proyecto" ---
$find4 = $d.Content.Find
$find4.ClearFormatting()
$find4.Replacement.ClearFormatting()
$find4.Execute("se procederá a ejecutar la garantía de fiel cumplimiento, de conformidad", $true, $false, $false, $false, $false, $true, 1, $false, `
    "se procederá a ejecutar las garantías asociadas al proyecto, de conformidad", 2) | Out-Null

$para4 = Get-ParagraphRange "*en tiempo y forma*"
Split-Sequence $para4.Start $para4.End @(
    "en tiempo y forma, situación que configura un incumplimiento reiterado, se procederá a ejecutar la",
    "s",
    " garantía",
    "s asociadas al proyecto"
)

# --- Edit 5: "Operaciones InnovaChile" -> split off "InnovaChile" -----
$para5 = Get-ParagraphRange "Operaciones InnovaChile"
Split-Sequence $para5.Start $para5.End @("Operaciones ", "InnovaChile")
